$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Change F1 header from "Change" to "Group"
$ws.Range("F1").Value = "Group"

# Update row 2 data with new Group value
$ws.Range("F2").Value = "A"

# Delete the old Before/After columns (G and H), shifting cells left
$ws.Range("G1:H3").Delete()
$ws.Range("F3").ClearContents()

# Apply the "Normal" style to the data rows (matches target formatting)
$ws.Range("A2:E3").Style = "Normal"

# Leave selection on the newly-entered Group value
$ws.Range("F2").Select() | Out-Null
